$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.329807877540588
$ws.Range("B1").Value = 1.905376553535461
$ws.Range("C1").Value = 1.845742702484131
$ws.Range("D1").Value = 4.693863391876221
$ws.Range("E1").Value = 1.300175905227661
